$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '41.294.71'
$ws.Range('E2').Value = '  -3.59%  '
$ws.Range('D3').Value = '2.462.28'
$ws.Range('E3').Value = '  -2.74%  '
$ws.Range('E4').Value = '  +0.03%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '312.02'
$ws.Range('D5').Style = $ws.Range('B5').Style
$ws.Range('E5').Value = '  +0.19%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '94.17'
$ws.Range('D6').Style = $ws.Range('B6').Style
$ws.Range('E6').Value = '  -6.89%  '
$ws.Range('E7').Value = '  -2.67%  '
$ws.Range('E8').Value = '  +0.04%  '
$ws.Range('E9').Value = '  -4.94%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '33.47'
$ws.Range('D10').Style = $ws.Range('B10').Style
$ws.Range('E10').Value = '  -6.71%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.0780'
$ws.Range('D11').Style = $ws.Range('B11').Style
$ws.Range('E11').Value = '  -3.33%  '
$ws.Range('E12').Value = '  -1.24%  '
$ws.Range('E13').Value = '  -4.87%  '
$ws.Range('D14').Value = '2.841.09'
$ws.Range('E14').Value = '  -2.80%  '
$ws.Range('D15').Value = '2.438.99'
$ws.Range('E15').Value = '  -4.35%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '14.84'
$ws.Range('D16').Style = $ws.Range('B16').Style
$ws.Range('E16').Value = '  -3.58%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '0.785'
$ws.Range('D17').Style = $ws.Range('B17').Style
$ws.Range('E17').Value = '  -3.94%  '
$ws.Range('D18').Value = '41.242.13'
$ws.Range('E18').Value = '  -3.67%  '
$ws.Range('E19').Value = '  -6.06%  '
$ws.Range('E20').Value = '  -3.36%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '11.27'
$ws.Range('D21').Style = $ws.Range('B21').Style
$ws.Range('E21').Value = '  -8.97%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '68.40'
$ws.Range('D22').Style = $ws.Range('B22').Style
$ws.Range('E22').Value = '  -1.88%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '236.81'
$ws.Range('D23').Style = $ws.Range('B23').Style
$ws.Range('E23').Value = '  -2.98%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '2.76'
$ws.Range('D24').Style = $ws.Range('B24').Style
$ws.Range('E24').Value = '  -4.14%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '1.90'
$ws.Range('D26').Style = $ws.Range('B26').Style
$ws.Range('E26').Value = '  -6.50%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '24.07'
$ws.Range('D27').Style = $ws.Range('B27').Style
$ws.Range('E27').Value = '  -5.59%  '
$ws.Range('E28').Value = '  -6.07%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '9.61'
$ws.Range('D29').Style = $ws.Range('B29').Style
$ws.Range('E29').Value = '  -5.74%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '36.65'
$ws.Range('D30').Style = $ws.Range('B30').Style
$ws.Range('E30').Value = '  -5.66%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '152.17'
$ws.Range('D31').Style = $ws.Range('B31').Style
$ws.Range('E31').Value = '  -6.09%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '5.49'
$ws.Range('D32').Style = $ws.Range('B32').Style
$ws.Range('E32').Value = '  -5.59%  '
$ws.Range('E33').Value = '  -5.20%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '2.59'
$ws.Range('D34').Style = $ws.Range('B34').Style
$ws.Range('E34').Value = '  -3.01%  '
$ws.Range('E35').Value = '  -5.87%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '3.06'
$ws.Range('D36').Style = $ws.Range('B36').Style
$ws.Range('E36').Value = '  -1.49%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '1.88'
$ws.Range('D37').Style = $ws.Range('B37').Style
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '16.97'
$ws.Range('D38').Style = $ws.Range('B38').Style
$ws.Range('E38').Value = '  -7.80%  '
$ws.Range('E39').Value = '  -3.11%  '
$ws.Range('E40').Value = '  -7.91%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '4.27'
$ws.Range('D41').Style = $ws.Range('B41').Style
$ws.Range('E41').Value = '  +1.59%  '
$ws.Range('E42').Value = '  +0.14%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '19.68'
$ws.Range('D43').Style = $ws.Range('B43').Style
$ws.Range('E43').Value = '  -11.56%  '
$ws.Range('D44').Value = '1.984.39'
$ws.Range('E44').Value = '  -0.18%  '
$ws.Range('E45').Value = '  -5.46%  '
$ws.Range('E46').Value = '  -8.36%  '
$ws.Range('E47').Value = '  -5.25%  '
$ws.Range('B48').Value = 'ordi'
$ws.Range('C48').Value = 'https://coinranking.com/coin/j7-7vPrOi+ordi-ordi'
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '68.99'
$ws.Range('D48').Style = $ws.Range('B48').Style
$ws.Range('E48').Value = '  -4.86%  '
$ws.Range('B49').Value = 'Aave'
$ws.Range('C49').Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '96.94'
$ws.Range('D49').Style = $ws.Range('B49').Style
$ws.Range('E49').Value = '  -4.07%  '
$ws.Range('E50').Value = '  -7.15%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '74.46'
$ws.Range('D51').Style = $ws.Range('B51').Style
$ws.Range('E51').Value = '  -6.47%  '
